$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused shared-formula cells in D/E for the blank separator
# rows (10, 17, 21) - they become plain empty (but still styled) cells.
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("D21").ClearContents()
$ws.Range("E21").ClearContents()

# Add the new task row 22: "create logo"
$ws.Range("A22").Value = "create logo"

$ws.Range("B22").NumberFormat = "0.00"
$ws.Range("B22").Value = 1

$ws.Range("C22").Value = 0.8

$ws.Range("D22").NumberFormat = "0.00"
$ws.Range("D22").Formula = "=C22"

$ws.Range("E22").NumberFormat = "0.00"
$ws.Range("E22").Formula = "=B22 - C22"

# Move the active selection like the author's session ended up (G24)
$ws.Range("G24").Select() | Out-Null
